# Update generated output numbers (view counts) across sheets.
$wb = $excel.ActiveWorkbook

$wsExhibition  = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAll         = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (Exhibition) ---
$wsExhibition.Range("F4").Value  = 8242
$wsExhibition.Range("F5").Value  = 6014
$wsExhibition.Range("F6").Value  = 514
$wsExhibition.Range("F7").Value  = 99
$wsExhibition.Range("F8").Value  = 17
$wsExhibition.Range("F11").Value = 880
$wsExhibition.Range("F12").Value = 79

# --- Sheet "演出" (Performance) ---
$wsPerformance.Range("F2").Value = 95

# --- Sheet "全部类型" (All types, combined view) ---
$wsAll.Range("F4").Value  = 8242
$wsAll.Range("F5").Value  = 6014
$wsAll.Range("F6").Value  = 514
$wsAll.Range("F7").Value  = 99
$wsAll.Range("F8").Value  = 17
$wsAll.Range("F11").Value = 95
$wsAll.Range("F15").Value = 880
$wsAll.Range("F16").Value = 79
